$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") rows 2-32: update the date serial value from 45212 to 45221
for ($r = 2; $r -le 32; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45212) {
        $cell.Value2 = 45221
    }
}
